$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column values are forced to text via a leading apostrophe (quote-prefix)
# so purely-numeric-looking strings (e.g. "1.00") are not coerced to numbers,
# matching the inline/shared string cells in the source workbook.
$priceCells = @()

$ws.Range("D2").Value = "'64.034.30"
$priceCells += "D2"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "'3.362.86"
$priceCells += "D3"
$ws.Range("E3").Value = "  +3.08%  "

$ws.Range("D4").Value = "'1.00"
$priceCells += "D4"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "'525.81"
$priceCells += "D5"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").Value = "'174.35"
$priceCells += "D6"
$ws.Range("E6").Value = "  -2.94%  "

$ws.Range("D7").Value = "'0.596"
$priceCells += "D7"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'3.360.16"
$priceCells += "D8"
$ws.Range("E8").Value = "  +3.02%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").Value = "'0.609"
$priceCells += "D10"
$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("D11").Value = "'53.37"
$priceCells += "D11"
$ws.Range("E11").Value = "  -7.97%  "

$ws.Range("D12").Value = "'0.134"
$priceCells += "D12"
$ws.Range("E12").Value = "  +2.68%  "

$ws.Range("D13").Value = "'0.0000257"
$priceCells += "D13"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").Value = "'9.09"
$priceCells += "D14"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "'3.897.60"
$priceCells += "D15"
$ws.Range("E15").Value = "  +3.59%  "

$ws.Range("D16").Value = "'3.360.34"
$priceCells += "D16"
$ws.Range("E16").Value = "  +3.44%  "

$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "'17.57"
$priceCells += "D18"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "'64.069.36"
$priceCells += "D19"
$ws.Range("E19").Value = "  +2.00%  "

$ws.Range("D20").Value = "'11.28"
$priceCells += "D20"
$ws.Range("E20").Value = "  +3.32%  "

$ws.Range("D21").Value = "'0.966"
$priceCells += "D21"
$ws.Range("E21").Value = "  +1.98%  "

$ws.Range("D22").Value = "'373.89"
$priceCells += "D22"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").Value = "'11.59"
$priceCells += "D23"
$ws.Range("E23").Value = "  +3.05%  "

$ws.Range("D24").Value = "'4.14"
$priceCells += "D24"
$ws.Range("E24").Value = "  +8.87%  "

$ws.Range("D25").Value = "'81.29"
$priceCells += "D25"
$ws.Range("E25").Value = "  +2.09%  "

$ws.Range("D26").Value = "'3.72"
$priceCells += "D26"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("D27").Value = "'6.15"
$priceCells += "D27"
$ws.Range("E27").Value = "  +1.76%  "

$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("D29").Value = "'11.32"
$priceCells += "D29"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").Value = "'8.25"
$priceCells += "D30"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("D31").Value = "'28.90"
$priceCells += "D31"
$ws.Range("E31").Value = "  +1.71%  "

$ws.Range("D32").Value = "'633.05"
$priceCells += "D32"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("D33").Value = "'6.44"
$priceCells += "D33"
$ws.Range("E33").Value = "  -4.57%  "

$ws.Range("D34").Value = "'11.22"
$priceCells += "D34"
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").Value = "'57.93"
$priceCells += "D36"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").Value = "'1.00"
$priceCells += "D37"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").Value = "'36.40"
$priceCells += "D38"
$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("D39").Value = "'0.380"
$priceCells += "D39"
$ws.Range("E39").Value = "  -4.37%  "

$ws.Range("D40").Value = "'0.0₃0743"
$priceCells += "D40"
$ws.Range("E40").Value = "  +13.81%  "

$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("D42").Value = "'2.66"
$priceCells += "D42"
$ws.Range("E42").Value = "  +8.80%  "

$ws.Range("D43").Value = "'2.980.57"
$priceCells += "D43"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("E44").Value = "  +0.93%  "

$ws.Range("D45").Value = "'2.97"
$priceCells += "D45"
$ws.Range("E45").Value = "  +5.90%  "

$ws.Range("D46").Value = "'2.69"
$priceCells += "D46"
$ws.Range("E46").Value = "  +4.30%  "

$ws.Range("D47").Value = "'0.0396"
$priceCells += "D47"
$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("D48").Value = "'2.61"
$priceCells += "D48"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("D49").Value = "'3.04"
$priceCells += "D49"
$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("D50").Value = "'0.125"
$priceCells += "D50"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("D51").Value = "'137.43"
$priceCells += "D51"
$ws.Range("E51").Value = "  +5.41%  "

# Reset number formatting back to the default "Normal" style so the
# quote-prefix formatting used above does not linger on the cell style.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
